# "Add alert when adding new data"
#
# The app now pops an alert whenever a new submission comes in; the sheet
# snapshot captured while testing that feature shows the extra rows that
# were appended while trying it out:
#   - Hoja1 (the form-responses sheet) gains rows 31-39
#   - SecondSheet (the grades sheet) gains rows 13-14
#
# The source workbook (built with SheetJS) stores every cell - even
# numeric-looking ones like DNI "442" - as text, never as a real number.
# We reproduce that by forcing a text number format ("@") and using the
# classic leading-apostrophe text prefix on the Formula, which keeps
# values such as "442" or "3376" from being auto-converted to numbers.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("SecondSheet")

function Set-TextCell($ws, $addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Formula = "'" + $value
}

# ---- Hoja1: new rows 31-39 (columns B..AB) ----
$sheet1Cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")
$sheet1NewRows = @{
    31 = @{ C = "442"; D = "qe"; E = "ff" }
    32 = @{ C = "3376" }
    33 = @{ C = "gre" }
    34 = @{}
    35 = @{}
    36 = @{}
    37 = @{}
    38 = @{}
    39 = @{ C = "42" }
}

foreach ($rowNum in ($sheet1NewRows.Keys | Sort-Object)) {
    $rowVals = $sheet1NewRows[$rowNum]
    foreach ($col in $sheet1Cols) {
        $val = ""
        if ($rowVals.ContainsKey($col)) { $val = $rowVals[$col] }
        $addr = $col + $rowNum
        Set-TextCell $ws1 $addr $val
    }
}

# ---- SecondSheet: new rows 13-14 (columns A,B,C,E,F,G,H - D stays blank/absent) ----
$sheet2Cols = @("A","B","C","E","F","G","H")
$sheet2NewRows = @{
    13 = @{ A = "3345214"; B = "";   C = "qqqqq"; E = "matematica"; F = "3"; G = "Desaprobado"; H = "4to" }
    14 = @{ A = "442";     B = "qe"; C = "ff";    E = "fr";         F = "5"; G = "Aprobado";    H = "wvwfwfqf" }
}

foreach ($rowNum in ($sheet2NewRows.Keys | Sort-Object)) {
    $rowVals = $sheet2NewRows[$rowNum]
    foreach ($col in $sheet2Cols) {
        $val = ""
        if ($rowVals.ContainsKey($col)) { $val = $rowVals[$col] }
        $addr = $col + $rowNum
        Set-TextCell $ws2 $addr $val
    }
}

Write-Host "Added Hoja1!B31:AB39 and SecondSheet!A13:H14"
